# Daily attendance processing - 2026-01-04 08:39:14
#
# Applies the daily attendance recalculation to the
# "Session Analysis Results" sheet:
#   - Summary counters (Missing/Pending sessions) refreshed
#   - "Recorded By" values re-ordered (System now listed first)
#   - Group-statistics Missing/Pending columns refreshed
#   - Sessions that rolled over to "Not Recorded" status get the
#     pink "Not Recorded" formatting (previously yellow "Pending")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary box (K/L columns near the top) ---------------------------
$ws.Range("L7").Value = 15   # Missing Sessions
$ws.Range("L8").Value = 96   # Pending Sessions

# --- "Recorded By" values: System now appears before the user --------
$recordedByRows = 8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309

foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# --- Group Statistics table: Missing / Pending columns refreshed -----
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 7

$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 7

$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 7

$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 7

$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 7

$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 7

# --- Sessions that moved from "Pending" to "Not Recorded" -------------
# Row 3 already carries the "Not Recorded" (pink) formatting - reuse it
# as the formatting template for the A:I block of each affected row.
$templateRange = $ws.Range("A3:I3")
$notRecordedRows = 20,46,72,98,124,150

foreach ($r in $notRecordedRows) {
    $templateRange.Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
    $ws.Range("I$r").Value = "Not Recorded"
}
